# AddressFinder.py now seeds the output with a starting address ("1 Rue du
# Clos Tellier") as package 0, re-running the geocoder and regenerating the
# whole "addresses found" table. The new run produced fewer rows (22 vs 31)
# with different addresses / geocoded coordinates / Maps links throughout.
# Apply the new table contents here and drop the rows that no longer exist.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 23-31 (old extra rows no longer present; final range is A1:E22)
$ws.Range("A23:E31").EntireRow.Delete()

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = '1 Rue du Clos Tellier'
$ws.Range("C2").Value = 49.377805
$ws.Range("D2").Value = 1.115311
$ws.Range("E2").Value = 'https://www.google.com/maps/place/1+Rue+du+Clos+Tellier'

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 'Rue Alfred Duthil,  La Girafe'
$ws.Range("C3").Value = 49.44899691233225
$ws.Range("D3").Value = 1.152946525663122
$ws.Range("E3").Value = 'https://www.google.com/maps/search/?api=1&query=49.44899691233225,1.1529465256631215'

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 'Rue du Four Banal,  Déville-lès-Rouen'
$ws.Range("C4").Value = 49.46194366399145
$ws.Range("D4").Value = 1.051755879188496
$ws.Range("E4").Value = 'https://www.google.com/maps/search/?api=1&query=49.46194366399145,1.0517558791884958'

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 'Centre hospitalier du Rouvray,  4'
$ws.Range("C5").Value = 49.39848725
$ws.Range("D5").Value = 1.09497675154378
$ws.Range("E5").Value = 'https://www.google.com/maps/search/?api=1&query=49.39848725,1.09497675154378'

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 'Avenue du Grand Cours,  La Sablière'
$ws.Range("C6").Value = 49.4260239
$ws.Range("D6").Value = 1.102107434450814
$ws.Range("E6").Value = 'https://www.google.com/maps/search/?api=1&query=49.426023900000004,1.1021074344508142'

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 'Boulevard Maurice de Broglie,  Parc de la Bresle'
$ws.Range("C7").Value = 49.46046790356267
$ws.Range("D7").Value = 1.073605849279293
$ws.Range("E7").Value = 'https://www.google.com/maps/search/?api=1&query=49.46046790356267,1.073605849279293'

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 'Foire Saint-Romain,  Presqu''Île Rollet'
$ws.Range("C8").Value = 49.44684465
$ws.Range("D8").Value = 1.054295428337912
$ws.Range("E8").Value = 'https://www.google.com/maps/search/?api=1&query=49.44684465,1.0542954283379116'

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 'Quai Richard Waddington,  Presqu''Île Rollet'
$ws.Range("C9").Value = 49.4407633
$ws.Range("D9").Value = 1.0466735
$ws.Range("E9").Value = 'https://www.google.com/maps/search/?api=1&query=49.4407633,1.0466735'

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 'Rue Michel,  Les Coquets'
$ws.Range("C10").Value = 49.45709401172554
$ws.Range("D10").Value = 1.091581960520784
$ws.Range("E10").Value = 'https://www.google.com/maps/search/?api=1&query=49.457094011725545,1.0915819605207844'

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 'Quai de Rouen-Quevilly,  Presqu''Île Rollet'
$ws.Range("C11").Value = 49.43784023857489
$ws.Range("D11").Value = 1.052916873215403
$ws.Range("E11").Value = 'https://www.google.com/maps/search/?api=1&query=49.437840238574886,1.052916873215403'

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 'Rue Ernest Lesueur,  Les Coquets'
$ws.Range("C12").Value = 49.45884816717249
$ws.Range("D12").Value = 1.085739919873663
$ws.Range("E12").Value = 'https://www.google.com/maps/search/?api=1&query=49.45884816717249,1.0857399198736635'

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 'Chemin de Rouen,  Le Mesnil-Esnard'
$ws.Range("C13").Value = 49.41840144619702
$ws.Range("D13").Value = 1.144385058694083
$ws.Range("E13").Value = 'https://www.google.com/maps/search/?api=1&query=49.418401446197024,1.1443850586940834'

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 'Rue François Lamy,  Mont-Riboudet'
$ws.Range("C14").Value = 49.4484727
$ws.Range("D14").Value = 1.0746356
$ws.Range("E14").Value = 'https://www.google.com/maps/search/?api=1&query=49.4484727,1.0746356'

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 'Hippodrome des Trois Pipes,  Rue de Diane'
$ws.Range("C15").Value = 49.45964395
$ws.Range("D15").Value = 1.118264972943661
$ws.Range("E15").Value = 'https://www.google.com/maps/search/?api=1&query=49.45964395,1.1182649729436607'

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 'Rue de Verdun,  Les Vikings'
$ws.Range("C16").Value = 49.47553723731352
$ws.Range("D16").Value = 1.123740339433279
$ws.Range("E16").Value = 'https://www.google.com/maps/search/?api=1&query=49.475537237313524,1.123740339433279'

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 'Chemin de la Bretèque,  La Bretèque'
$ws.Range("C17").Value = 49.4895594
$ws.Range("D17").Value = 1.1035399
$ws.Range("E17").Value = 'https://www.google.com/maps/search/?api=1&query=49.4895594,1.1035399'

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = '6,  Rue des Canadiens'
$ws.Range("C18").Value = 49.420784
$ws.Range("D18").Value = 1.128584978738117
$ws.Range("E18").Value = 'https://www.google.com/maps/search/?api=1&query=49.420784,1.1285849787381173'

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = '2,  Rue Legouy'
$ws.Range("C19").Value = 49.4425295
$ws.Range("D19").Value = 1.108076
$ws.Range("E19").Value = 'https://www.google.com/maps/search/?api=1&query=49.4425295,1.108076'

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 'Avenue Bernard Bicheray,  Rouen'
$ws.Range("C20").Value = 49.44945537954387
$ws.Range("D20").Value = 1.049347084850659
$ws.Range("E20").Value = 'https://www.google.com/maps/search/?api=1&query=49.44945537954387,1.0493470848506585'

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 'Rue Gessard,  Saint-Clément'
$ws.Range("C21").Value = 49.42722088104083
$ws.Range("D21").Value = 1.071434974292933
$ws.Range("E21").Value = 'https://www.google.com/maps/search/?api=1&query=49.42722088104083,1.0714349742929332'

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 'Rampe Beauvoisine,  Jouvenet'
$ws.Range("C22").Value = 49.44801269803951
$ws.Range("D22").Value = 1.102314470687749
$ws.Range("E22").Value = 'https://www.google.com/maps/search/?api=1&query=49.44801269803951,1.102314470687749'

